$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style/format from N1 into O1:Q1 (matches existing bold/border/center style s="1")
$ws.Range("N1").Copy()
$ws.Range("O1:Q1").PasteSpecial(-4122)
$ws.Range("O1").Value = "31/12/2023"
$ws.Range("P1").Value = "31/03/2024"
$ws.Range("Q1").Value = "30/06/2024"

# Data rows 2-80: columns O (15), P (16), Q (17)
$newData = @{
    2 = @(20119857.152, 20232347.648, 21740724.224)
    3 = @(3965496.064, 3556328.96, 4119134.976)
    4 = @(1754105.984, 1019854.976, 1533334.016)
    5 = @(154559.008, 9123, 50630)
    6 = @(522022.016, 496428.992, 459680.992)
    7 = @(814819.008, 937560, 799728)
    8 = @(0, 0, 0)
    9 = @(159898, 277356, 370311.008)
    10 = @(164556, 173338, 144159.008)
    11 = @(395536, 642668.032, 761292.032)
    12 = @(2940807.936, 3207396.096, 3759728.896)
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(2304150.016, 2498099.968, 2779450.112)
    16 = @(0, 0, 0)
    17 = @(0, 0, 0)
    18 = @(0, 0, 0)
    19 = @(0, 0, 0)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(1600, 1600, 1600)
    23 = @(6149094.912, 6504949.248, 7036836.864)
    24 = @(7062858.752, 6962073.088, 6823421.952)
    25 = @(0, 0, 0)
    26 = @(20119857.152, 20232347.648, 21740724.224)
    27 = @(3336510.976, 2807996.928, 3010781.952)
    28 = @(103832, 81940, 102487)
    29 = @(1315213.952, 1208488.96, 1299776)
    30 = @(29376, 73543, 133918)
    31 = @(961353.024, 238903.008, 309888.992)
    32 = @(28293, 28154, 37408)
    33 = @(92565, 92565, 92565)
    34 = @(766985.0159999999, 1049496.952, 993025.0160000001)
    35 = @(38893, 34906, 41714)
    36 = @(0, 0, 0)
    37 = @(11218256.896, 12067040.256, 13705895.936)
    38 = @(8300813.824, 9331755.007999999, 10802684.928)
    39 = @(16071, 10718, 5357)
    40 = @(1480519.04, 1255885.952, 1376204.032)
    41 = @(68288, 66739, 65103)
    42 = @(0, 0, 0)
    43 = @(1352564.992, 1401942.016, 1456546.944)
    44 = @(0, 0, 0)
    45 = @(0, 0, 0)
    46 = @(72653, 78390, 79980)
    47 = @(5492436.792, 5278919.952, 4944065.056)
    48 = @(5055782.912, 5062063.104, 5062422.016)
    49 = @(139471.008, 155296.992, 184728)
    50 = @(0, 0, 0)
    51 = @(297183.008, 297183.008, 297183.008)
    52 = @(0, -235623.008, -600268.032)
    53 = @(0, 0, 0)
    54 = @(0, 0, 0)
    55 = @(0, 0, 0)
    56 = @(0, 0, 0)
    59 = @(1849449.984, 2007601.024, 2575361.024)
    60 = @(-1270588.032, -1356478.976, -1843384.96)
    61 = @(578861.952, 651121.9840000001, 731976)
    62 = @(0, 0, 0)
    63 = @(-90168, -136474, -133834)
    64 = @(26500, 0, 0)
    65 = @(-48130, 0, 0)
    66 = @(44030, -8182, -18250)
    67 = @(0, 0, 0)
    68 = @(36761.032, -765401.024, -1125346.944)
    69 = @(219843.984, 117388, 208463.008)
    70 = @(-183082.992, -882788.992, -1333810.048)
    74 = @(547854.976, -258935.008, -545454.976)
    75 = @(-53973, -73537, -60373)
    76 = @(-86704.992, 102586, 242772.992)
    79 = @(-8162, -5737, -1590)
    80 = @(399014.976, -235623.008, -365644.992)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 15).Value = $vals[0]
    $ws.Cells.Item($r, 16).Value = $vals[1]
    $ws.Cells.Item($r, 17).Value = $vals[2]
}

Write-Output "Applied edits"